# positions.xlsx — add grant rows, adjust content
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Grants section (rows 26-28): replace old PI grants with new, detailed grants ---
$ws.Range("A26").Value = "博士后基金面上项目"
$ws.Range("B26").Value = " CNY ¥5（主持）"
$ws.Range("C26").Value = "NA"
$ws.Range("D26").Value = "2013"
$ws.Range("E26").Value = "2016"
$ws.Range("F26").Value = "NA"
$ws.Range("G26").Value = "NA"
$ws.Range("H26").Value = "NA"
$ws.Range("I26").Value = "grant"
$ws.Range("J26").Value = "NA"

$ws.Range("A27").Value = "国家重点研发计划子课题"
$ws.Range("B27").Value = " CNY ¥72.5（主持）"
$ws.Range("C27").Value = "NA"
$ws.Range("D27").Value = "2016"
$ws.Range("E27").Value = "2021"
$ws.Range("F27").Value = "NA"
$ws.Range("G27").Value = "NA"
$ws.Range("H27").Value = "NA"
$ws.Range("I27").Value = "grant"
$ws.Range("J27").Value = "NA"

$ws.Range("A28").Value = "NSFC 青年科学基金"
$ws.Range("B28").Value = " CNY ¥24（主持）"
$ws.Range("C28").Value = "NA"
$ws.Range("D28").Value = "2022"
$ws.Range("E28").Value = "2024"
$ws.Range("F28").Value = "NA"
$ws.Range("G28").Value = "NA"
$ws.Range("H28").Value = "NA"
$ws.Range("I28").Value = "grant"
$ws.Range("J28").Value = "NA"

# --- Conference presentation row 32: mark loc with leading "#" ---
$ws.Range("A32").Value = "# 华中农业大学资源与环境学院博士后交流会"

# --- Conferences row 43: toggle "_present_" to "# _present_" ---
$ws.Range("A43").Value = "# _present_"

# --- Update the active selection to A38 ---
$ws.Range("A38").Select()
